# Scheduled market-data refresh for Pandaemonium Profits leve tables.
# Updates currentAveragePrice* / LevePrice* / LeveProfit* columns (H:N)
# per sheet/row with freshly pulled Universalis price data.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 398
$ws.Range("I8").Value = 398
$ws.Range("K8").Value = 1194
$ws.Range("M8").Value = -1055

$ws.Range("H40").Value = 2291.3914
$ws.Range("I40").Value = 2744.9092
$ws.Range("J40").Value = 1875.6666
$ws.Range("K40").Value = 2744.9092
$ws.Range("L40").Value = 1875.6666
$ws.Range("M40").Value = -2569.9092
$ws.Range("N40").Value = -2225.6666

$ws.Range("H70").Value = 1541.7646
$ws.Range("I70").Value = 1408.6
$ws.Range("J70").Value = 1732
$ws.Range("K70").Value = 4225.799999999999
$ws.Range("L70").Value = 5196
$ws.Range("M70").Value = -3955.799999999999
$ws.Range("N70").Value = -5736

$ws.Range("H73").Value = 1541.7646
$ws.Range("I73").Value = 1408.6
$ws.Range("J73").Value = 1732
$ws.Range("K73").Value = 4225.799999999999
$ws.Range("L73").Value = 5196
$ws.Range("M73").Value = -3289.799999999999
$ws.Range("N73").Value = -7068

$ws.Range("H76").Value = 3599.8684
$ws.Range("I76").Value = 3503.1667
$ws.Range("K76").Value = 3503.1667
$ws.Range("M76").Value = -3188.1667

$ws.Range("H79").Value = 3599.8684
$ws.Range("I79").Value = 3503.1667
$ws.Range("K79").Value = 3503.1667
$ws.Range("M79").Value = -2411.1667

$ws.Range("H107").Value = 1321
$ws.Range("I107").Value = 1201.25
$ws.Range("J107").Value = 1800
$ws.Range("K107").Value = 1201.25
$ws.Range("L107").Value = 1800
$ws.Range("M107").Value = 718.75
$ws.Range("N107").Value = -5640

$ws.Range("H129").Value = 1032.5128
$ws.Range("J129").Value = 1032.5128
$ws.Range("L129").Value = 3097.5384
$ws.Range("N129").Value = -13097.5384

$ws.Range("H132").Value = 1942.814
$ws.Range("I132").Value = 1595.5834
$ws.Range("J132").Value = 3728.5715
$ws.Range("K132").Value = 4786.7502
$ws.Range("L132").Value = 11185.7145
$ws.Range("M132").Value = -2256.7502
$ws.Range("N132").Value = -16245.7145

$ws.Range("H137").Value = 593494.1
$ws.Range("I137").Value = 1573.3793
$ws.Range("J137").Value = 1547144.1
$ws.Range("K137").Value = 4720.1379
$ws.Range("L137").Value = 4641432.300000001
$ws.Range("M137").Value = -2170.1379
$ws.Range("N137").Value = -4646532.300000001

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H63").Value = 2500
$ws.Range("I63").Value = 2500
$ws.Range("K63").Value = 2500
$ws.Range("M63").Value = -1814

$ws.Range("H66").Value = 2500
$ws.Range("I66").Value = 2500
$ws.Range("K66").Value = 12500
$ws.Range("M66").Value = -9068

$ws.Range("H97").Value = 585.2857
$ws.Range("I97").Value = 585.2857
$ws.Range("K97").Value = 585.2857
$ws.Range("M97").Value = -89.28570000000002

$ws.Range("H107").Value = 39228
$ws.Range("J107").Value = 39228
$ws.Range("L107").Value = 39228
$ws.Range("N107").Value = -46908

$ws.Range("H110").Value = 0
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 0
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 0
$ws.Range("M110").Value = ""
$ws.Range("N110").Value = ""

$ws.Range("H111").Value = 79800
$ws.Range("J111").Value = 79800
$ws.Range("L111").Value = 79800
$ws.Range("N111").Value = -87980

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1904.3829
$ws.Range("I86").Value = 1964.3334
$ws.Range("J86").Value = 1612.125
$ws.Range("K86").Value = 1964.3334
$ws.Range("L86").Value = 1612.125
$ws.Range("M86").Value = -841.3334
$ws.Range("N86").Value = -3858.125

$ws.Range("H89").Value = 1904.3829
$ws.Range("I89").Value = 1964.3334
$ws.Range("J89").Value = 1612.125
$ws.Range("K89").Value = 9821.666999999999
$ws.Range("L89").Value = 8060.625
$ws.Range("M89").Value = -4205.666999999999
$ws.Range("N89").Value = -19292.625

$ws.Range("H99").Value = 1665.5
$ws.Range("I99").Value = 1473.579
$ws.Range("J99").Value = 1946
$ws.Range("K99").Value = 1473.579
$ws.Range("L99").Value = 1946
$ws.Range("M99").Value = 24.42100000000005
$ws.Range("N99").Value = -4942

$ws.Range("H134").Value = 5325.2085
$ws.Range("I134").Value = 4611.2354
$ws.Range("J134").Value = 7059.143
$ws.Range("K134").Value = 13833.7062
$ws.Range("L134").Value = 21177.429
$ws.Range("M134").Value = -11298.7062
$ws.Range("N134").Value = -26247.429

$ws.Range("H138").Value = 73313.336
$ws.Range("J138").Value = 73313.336
$ws.Range("L138").Value = 73313.336
$ws.Range("N138").Value = -83593.336

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 140
$ws.Range("I22").Value = 103.333336
$ws.Range("J22").Value = 250
$ws.Range("K22").Value = 103.333336
$ws.Range("L22").Value = 250
$ws.Range("M22").Value = 246.666664
$ws.Range("N22").Value = -950

$ws.Range("H31").Value = 2538.0476
$ws.Range("I31").Value = 1812
$ws.Range("J31").Value = 4861.4
$ws.Range("K31").Value = 1812
$ws.Range("L31").Value = 4861.4
$ws.Range("M31").Value = -1517
$ws.Range("N31").Value = -5451.4

$ws.Range("H34").Value = 2538.0476
$ws.Range("I34").Value = 1812
$ws.Range("J34").Value = 4861.4
$ws.Range("K34").Value = 1812
$ws.Range("L34").Value = 4861.4
$ws.Range("M34").Value = -1610
$ws.Range("N34").Value = -5265.4

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H17").Value = 75.333336
$ws.Range("I17").Value = 75.333336
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 226.000008
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -57.00000800000001
$ws.Range("N17").Value = ""

$ws.Range("H32").Value = 1737.875
$ws.Range("J32").Value = 2067.1667
$ws.Range("L32").Value = 6201.500100000001
$ws.Range("N32").Value = -6767.500100000001

$ws.Range("H70").Value = 4316.6665
$ws.Range("I70").Value = 1000
$ws.Range("J70").Value = 4980
$ws.Range("K70").Value = 3000
$ws.Range("L70").Value = 14940
$ws.Range("M70").Value = -2685
$ws.Range("N70").Value = -15570

$ws.Range("H73").Value = 4316.6665
$ws.Range("I73").Value = 1000
$ws.Range("J73").Value = 4980
$ws.Range("K73").Value = 3000
$ws.Range("L73").Value = 14940
$ws.Range("M73").Value = -1908
$ws.Range("N73").Value = -17124

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 6920.643
$ws.Range("I80").Value = 18166.334
$ws.Range("K80").Value = 18166.334
$ws.Range("M80").Value = -17168.334

$ws.Range("H83").Value = 6920.643
$ws.Range("I83").Value = 18166.334
$ws.Range("K83").Value = 90831.67
$ws.Range("M83").Value = -85839.67

$ws.Range("H140").Value = 57773.2
$ws.Range("J140").Value = 57773.2
$ws.Range("L140").Value = 57773.2
$ws.Range("N140").Value = -68133.2

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 3701.9429
$ws.Range("I40").Value = 3704.6956
$ws.Range("J40").Value = 3696.6667
$ws.Range("K40").Value = 3704.6956
$ws.Range("L40").Value = 3696.6667
$ws.Range("M40").Value = -3568.6956
$ws.Range("N40").Value = -3968.6667

$ws.Range("H132").Value = 3747.3845
$ws.Range("I132").Value = 3189.276
$ws.Range("J132").Value = 5365.9
$ws.Range("K132").Value = 9567.828
$ws.Range("L132").Value = 16097.7
$ws.Range("M132").Value = -7037.828
$ws.Range("N132").Value = -21157.7

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4120.3
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 4300.75
$ws.Range("K62").Value = 4000
$ws.Range("L62").Value = 4300.75
$ws.Range("M62").Value = -3376
$ws.Range("N62").Value = -5548.75

$ws.Range("H65").Value = 4120.3
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 4300.75
$ws.Range("K65").Value = 20000
$ws.Range("L65").Value = 21503.75
$ws.Range("M65").Value = -16880
$ws.Range("N65").Value = -27743.75

$ws.Range("H132").Value = 1330.125
$ws.Range("I132").Value = 1574.25
$ws.Range("J132").Value = 963.9375
$ws.Range("K132").Value = 4722.75
$ws.Range("L132").Value = 2891.8125
$ws.Range("M132").Value = -2192.75
$ws.Range("N132").Value = -7951.8125
